$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 421.85715
$ws.Range("I39").Value = 221.6
$ws.Range("J39").Value = 501.96
$ws.Range("K39").Value = 664.8
$ws.Range("L39").Value = 1505.88
$ws.Range("M39").Value = -368.8
$ws.Range("N39").Value = -2097.88
$ws.Range("H40").Value = 3653.3428
$ws.Range("I40").Value = 3373.3333
$ws.Range("J40").Value = 3863.35
$ws.Range("K40").Value = 3373.3333
$ws.Range("L40").Value = 3863.35
$ws.Range("M40").Value = -3198.3333
$ws.Range("N40").Value = -4213.35
$ws.Range("H48").Value = 1466.6666
$ws.Range("I48").Value = 400
$ws.Range("K48").Value = 1200
$ws.Range("M48").Value = -908
$ws.Range("H56").Value = 1466.6666
$ws.Range("I56").Value = 400
$ws.Range("K56").Value = 1200
$ws.Range("M56").Value = -666
$ws.Range("H98").Value = 52633050
$ws.Range("I98").Value = 52633050
$ws.Range("K98").Value = 52633050
$ws.Range("M98").Value = -52631552
$ws.Range("H122").Value = 52633050
$ws.Range("I122").Value = 52633050
$ws.Range("K122").Value = 157899150
$ws.Range("M122").Value = -157896700
$ws.Range("H132").Value = 1393.3077
$ws.Range("I132").Value = 1393.3077
$ws.Range("K132").Value = 4179.9231
$ws.Range("M132").Value = -1649.9231
$ws.Range("H138").Value = 1718.5238
$ws.Range("I138").Value = 842.2857
$ws.Range("J138").Value = 3471
$ws.Range("K138").Value = 2526.8571
$ws.Range("L138").Value = 10413
$ws.Range("M138").Value = 2613.1429
$ws.Range("N138").Value = -20693
$ws.Range("H141").Value = 1961.8
$ws.Range("I141").Value = 1926
$ws.Range("K141").Value = 5778
$ws.Range("M141").Value = -598

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 12001
$ws.Range("I38").Value = 12001
$ws.Range("K38").Value = 12001
$ws.Range("M38").Value = -11534
$ws.Range("H45").Value = 2159.7917
$ws.Range("I45").Value = 1851.7
$ws.Range("J45").Value = 2379.8572
$ws.Range("K45").Value = 1851.7
$ws.Range("L45").Value = 2379.8572
$ws.Range("M45").Value = -1474.7
$ws.Range("N45").Value = -3133.8572
$ws.Range("H61").Value = 38544740
$ws.Range("I61").Value = 100002800
$ws.Range("K61").Value = 100002800
$ws.Range("M61").Value = -100002588
$ws.Range("H110").Value = 4041.3333
$ws.Range("I110").Value = 3511
$ws.Range("K110").Value = 3511
$ws.Range("M110").Value = -1466
$ws.Range("H136").Value = 38544740
$ws.Range("I136").Value = 100002800
$ws.Range("K136").Value = 300008400
$ws.Range("M136").Value = -300005850

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5199.1
$ws.Range("I20").Value = 5298
$ws.Range("K20").Value = 5298
$ws.Range("M20").Value = -5051
$ws.Range("H26").Value = 38000
$ws.Range("I26").Value = 38000
$ws.Range("K26").Value = 38000
$ws.Range("M26").Value = -37708
$ws.Range("H94").Value = 1122.8334
$ws.Range("I94").Value = 1140.4286
$ws.Range("J94").Value = 999.6667
$ws.Range("K94").Value = 1140.4286
$ws.Range("L94").Value = 999.6667
$ws.Range("M94").Value = -689.4286
$ws.Range("N94").Value = -1901.6667
$ws.Range("H96").Value = 33459
$ws.Range("J96").Value = 62557.668
$ws.Range("L96").Value = 62557.668
$ws.Range("N96").Value = -68049.66800000001
$ws.Range("H107").Value = 2400.875
$ws.Range("I107").Value = 2400.875
$ws.Range("K107").Value = 2400.875
$ws.Range("M107").Value = -480.875
$ws.Range("H134").Value = 32919.332
$ws.Range("I134").Value = 5011.9688
$ws.Range("K134").Value = 15035.9064
$ws.Range("M134").Value = -12500.9064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2806.8
$ws.Range("I16").Value = 2342.6667
$ws.Range("J16").Value = 3503
$ws.Range("K16").Value = 2342.6667
$ws.Range("L16").Value = 3503
$ws.Range("M16").Value = -2055.6667
$ws.Range("N16").Value = -4077
$ws.Range("H113").Value = 2806.8
$ws.Range("I113").Value = 2342.6667
$ws.Range("J113").Value = 3503
$ws.Range("K113").Value = 2342.6667
$ws.Range("L113").Value = 3503
$ws.Range("M113").Value = -172.6667000000002
$ws.Range("N113").Value = -7843
$ws.Range("H132").Value = 1713.9474
$ws.Range("I132").Value = 1289.875
$ws.Range("J132").Value = 3975.6667
$ws.Range("K132").Value = 3869.625
$ws.Range("L132").Value = 11927.0001
$ws.Range("M132").Value = -1339.625
$ws.Range("N132").Value = -16987.0001
$ws.Range("H134").Value = 670566.25
$ws.Range("I134").Value = 1112720.9
$ws.Range("K134").Value = 3338162.7
$ws.Range("M134").Value = -3335627.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 13983.333
$ws.Range("I3").Value = 1950
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 5850
$ws.Range("L3").Value = 60000
$ws.Range("M3").Value = -5738
$ws.Range("N3").Value = -60224
$ws.Range("H39").Value = 107820.93
$ws.Range("J39").Value = 105799.9
$ws.Range("L39").Value = 317399.7
$ws.Range("N39").Value = -317987.7
$ws.Range("H129").Value = 83581180
$ws.Range("I129").Value = 1930
$ws.Range("K129").Value = 5790
$ws.Range("M129").Value = -790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H93").Value = 59995
$ws.Range("J93").Value = 59995
$ws.Range("L93").Value = 59995
$ws.Range("N93").Value = -63739
$ws.Range("H126").Value = 6001.75
$ws.Range("I126").Value = 7250
$ws.Range("K126").Value = 21750
$ws.Range("M126").Value = -19280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4993.75
$ws.Range("I16").Value = 4991.6665
$ws.Range("K16").Value = 4991.6665
$ws.Range("M16").Value = -4821.6665
$ws.Range("H68").Value = 4260.6
$ws.Range("I68").Value = 1825.75
$ws.Range("K68").Value = 1825.75
$ws.Range("M68").Value = -1076.75
$ws.Range("H71").Value = 4260.6
$ws.Range("I71").Value = 1825.75
$ws.Range("K71").Value = 9128.75
$ws.Range("M71").Value = -5384.75
$ws.Range("H82").Value = 3000.3333
$ws.Range("I82").Value = 3002
$ws.Range("J82").Value = 2999.5
$ws.Range("K82").Value = 3002
$ws.Range("L82").Value = 2999.5
$ws.Range("M82").Value = -2641
$ws.Range("N82").Value = -3721.5
$ws.Range("H85").Value = 3000.3333
$ws.Range("I85").Value = 3002
$ws.Range("J85").Value = 2999.5
$ws.Range("K85").Value = 3002
$ws.Range("L85").Value = 2999.5
$ws.Range("M85").Value = -1754
$ws.Range("N85").Value = -5495.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 14752.25
$ws.Range("J11").Value = 14752.25
$ws.Range("L11").Value = 14752.25
$ws.Range("N11").Value = -15036.25
$ws.Range("H15").Value = 36335.668
$ws.Range("I15").Value = 20000
$ws.Range("J15").Value = 44503.5
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 44503.5
$ws.Range("M15").Value = -19712
$ws.Range("N15").Value = -45079.5
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("I21").Value = 75000
$ws.Range("K21").Value = 75000
$ws.Range("M21").Value = -74765
$ws.Range("I35").Value = 75000
$ws.Range("K35").Value = 75000
$ws.Range("M35").Value = -74710
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H52").Value = 56845.668
$ws.Range("I52").Value = 65021
$ws.Range("K52").Value = 65021
$ws.Range("M52").Value = -64795
$ws.Range("H136").Value = 1486.125
$ws.Range("I136").Value = 1404.7368
$ws.Range("K136").Value = 4214.2104
$ws.Range("M136").Value = -1664.2104
